# "add changelog, highlight md links, update config about_main"
#
# config.xlsx / Sheet1 is a 2-column id/value settings table (A=id, B=value).
# Row 3 holds the "about_main" markdown blurb shown on the app's home page.
# Replace it with the refreshed copy: dropped bold-markdown emphasis on a few
# phrases, dropped the "Indépendant" bullet, reworded "Portable"/"Sécurisé",
# added "Doc" to the concept list, and added a link to datannur.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newAboutMain = @"
![main_banner not_rounded](data/img/main_banner{dark_mode}.png?v=1)

# datannur, le catalogue de données portable

Permet de centraliser, rechercher et visualiser les informations sur une collection de jeux de données

Pour améliorer l’organisation des données et faciliter leur partage et leur documentation

Simple et flexible, s’intègre rapidement dans tous types d’environnement


- **Facile** :
Aucune installation ou configuration nécessaire. datannur est le catalogue le plus simple à implémenter et maintenir

- **Portable** :
Fonctionne partout (local, cloud, disque partagé), un dossier que l’on peut copier, déplacer, envoyer et ouvrir avec n’importe quel navigateur

- **Complet** :
Flexible, complet et structuré autour de 7 concepts avec un niveau de détail important : Institution, Dossier, Mot clé, Doc, Dataset, Variable et Modalité

- **Sécurisé** :
Parce qu’elle est une simple interface HTML isolée dans le navigateur, l’application ne peut rien modifier sur la machine et ne pose ainsi aucun risque

Pour davantage d'information : [datannur.com](https://datannur.com)

La version ici présente est un prototype en cours de développement et d'expérimentation. Les données utilisées sont fictives et uniquement à usage de test et de développement. Question ou suggestion : [contact@datannur.com](mailto:contact@datannur.com).
"@

$ws.Range("B3").Value = $newAboutMain

# Leave the cursor on the cell that was just edited (was B22).
$ws.Range("B3").Select()
